$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2270.9524
$ws.Range("J17").Value = 2270.9524
$ws.Range("L17").Value = 6812.8572
$ws.Range("N17").Value = -7148.8572

$ws.Range("H64").Value = 166670420
$ws.Range("I64").Value = 1000000000
$ws.Range("J64").Value = 4492
$ws.Range("K64").Value = 1000000000
$ws.Range("L64").Value = 4492
$ws.Range("M64").Value = -999999752
$ws.Range("N64").Value = -4988

$ws.Range("H67").Value = 166670420
$ws.Range("I67").Value = 1000000000
$ws.Range("J67").Value = 4492
$ws.Range("K67").Value = 1000000000
$ws.Range("L67").Value = 4492
$ws.Range("M67").Value = -999999142
$ws.Range("N67").Value = -6208

$ws.Range("H74").Value = 2420.149
$ws.Range("I74").Value = 2030.1923
$ws.Range("J74").Value = 2902.9524
$ws.Range("K74").Value = 2030.1923
$ws.Range("L74").Value = 2902.9524
$ws.Range("M74").Value = -1094.1923
$ws.Range("N74").Value = -4774.9524

$ws.Range("H75").Value = 31416.5
$ws.Range("J75").Value = 31416.5
$ws.Range("L75").Value = 31416.5
$ws.Range("N75").Value = -33288.5

$ws.Range("H77").Value = 2420.149
$ws.Range("I77").Value = 2030.1923
$ws.Range("J77").Value = 2902.9524
$ws.Range("K77").Value = 10150.9615
$ws.Range("L77").Value = 14514.762
$ws.Range("M77").Value = -5470.961499999999
$ws.Range("N77").Value = -23874.762

$ws.Range("H78").Value = 31416.5
$ws.Range("J78").Value = 31416.5
$ws.Range("L78").Value = 94249.5
$ws.Range("N78").Value = -103609.5

$ws.Range("H107").Value = 201.53847
$ws.Range("I107").Value = 185.45454
$ws.Range("J107").Value = 290
$ws.Range("K107").Value = 185.45454
$ws.Range("L107").Value = 290
$ws.Range("M107").Value = 1734.54546
$ws.Range("N107").Value = -4130

$ws.Range("H116").Value = 3087
$ws.Range("I116").Value = 2607.5293
$ws.Range("K116").Value = 2607.5293
$ws.Range("M116").Value = 834.4706999999999

$ws.Range("H118").Value = 2058.2334
$ws.Range("I118").Value = 961.75
$ws.Range("J118").Value = 2226.923
$ws.Range("K118").Value = 2885.25
$ws.Range("L118").Value = 6680.768999999999
$ws.Range("M118").Value = -1228.25
$ws.Range("N118").Value = -9994.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 260703.62
$ws.Range("I32").Value = 269300.2
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 269300.2
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -269013.2
$ws.Range("N32").Value = -20574

$ws.Range("H63").Value = 100
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 100
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 100
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -1472

$ws.Range("H66").Value = 100
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 100
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 500
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -7364

$ws.Range("H74").Value = 5503.8696
$ws.Range("I74").Value = 1000.8823
$ws.Range("J74").Value = 18262.334
$ws.Range("K74").Value = 1000.8823
$ws.Range("L74").Value = 18262.334
$ws.Range("M74").Value = -126.8823
$ws.Range("N74").Value = -20010.334

$ws.Range("H77").Value = 5503.8696
$ws.Range("I77").Value = 1000.8823
$ws.Range("J77").Value = 18262.334
$ws.Range("K77").Value = 5004.4115
$ws.Range("L77").Value = 91311.67
$ws.Range("M77").Value = -636.4115000000002
$ws.Range("N77").Value = -100047.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 515.5
$ws.Range("I22").Value = 245.375
$ws.Range("J22").Value = 2676.5
$ws.Range("K22").Value = 245.375
$ws.Range("L22").Value = 2676.5
$ws.Range("M22").Value = 104.625
$ws.Range("N22").Value = -3376.5

$ws.Range("H31").Value = 14203.117
$ws.Range("I31").Value = 1458.5483
$ws.Range("J31").Value = 20474.254
$ws.Range("K31").Value = 1458.5483
$ws.Range("L31").Value = 20474.254
$ws.Range("M31").Value = -1163.5483
$ws.Range("N31").Value = -21064.254

$ws.Range("H34").Value = 14203.117
$ws.Range("I34").Value = 1458.5483
$ws.Range("J34").Value = 20474.254
$ws.Range("K34").Value = 1458.5483
$ws.Range("L34").Value = 20474.254
$ws.Range("M34").Value = -1256.5483
$ws.Range("N34").Value = -20878.254

$ws.Range("H58").Value = 962.7931
$ws.Range("I58").Value = 538.25
$ws.Range("J58").Value = 1906.2222
$ws.Range("K58").Value = 538.25
$ws.Range("L58").Value = 1906.2222
$ws.Range("M58").Value = -335.25
$ws.Range("N58").Value = -2312.2222

$ws.Range("H68").Value = 20850
$ws.Range("J68").Value = 23300
$ws.Range("L68").Value = 23300
$ws.Range("N68").Value = -24798

$ws.Range("H71").Value = 20850
$ws.Range("J71").Value = 23300
$ws.Range("L71").Value = 69900
$ws.Range("N71").Value = -77388

$ws.Range("H74").Value = 20314
$ws.Range("J74").Value = 20314
$ws.Range("L74").Value = 20314
$ws.Range("N74").Value = -22062

$ws.Range("H77").Value = 20314
$ws.Range("J77").Value = 20314
$ws.Range("L77").Value = 60942
$ws.Range("N77").Value = -69678

$ws.Range("H87").Value = 17499.75
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 17499.75
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 17499.75
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -19871.75

$ws.Range("H90").Value = 17499.75
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 17499.75
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 52499.25
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -64355.25

$ws.Range("H132").Value = 27031780
$ws.Range("I132").Value = 31254900
$ws.Range("J132").Value = 3808
$ws.Range("K132").Value = 93764700
$ws.Range("L132").Value = 11424
$ws.Range("M132").Value = -93762170
$ws.Range("N132").Value = -16484

$ws.Range("H136").Value = 962.7931
$ws.Range("I136").Value = 538.25
$ws.Range("J136").Value = 1906.2222
$ws.Range("K136").Value = 1614.75
$ws.Range("L136").Value = 5718.6666
$ws.Range("M136").Value = 935.25
$ws.Range("N136").Value = -10818.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4863647
$ws.Range("I137").Value = 93464.164
$ws.Range("J137").Value = 11223891
$ws.Range("K137").Value = 280392.492
$ws.Range("L137").Value = 33671673
$ws.Range("M137").Value = -275292.492
$ws.Range("N137").Value = -33681873

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3041.7144
$ws.Range("I80").Value = 2124
$ws.Range("J80").Value = 3551.5557
$ws.Range("K80").Value = 2124
$ws.Range("L80").Value = 3551.5557
$ws.Range("M80").Value = -1126
$ws.Range("N80").Value = -5547.5557

$ws.Range("H83").Value = 3041.7144
$ws.Range("I83").Value = 2124
$ws.Range("J83").Value = 3551.5557
$ws.Range("K83").Value = 10620
$ws.Range("L83").Value = 17757.7785
$ws.Range("M83").Value = -5628
$ws.Range("N83").Value = -27741.7785

$ws.Range("H122").Value = 928.9167
$ws.Range("I122").Value = 899.7273
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 2699.1819
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -249.1819
$ws.Range("N122").Value = -8650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1347.7
$ws.Range("I82").Value = 1330.6666
$ws.Range("J82").Value = 1355
$ws.Range("K82").Value = 1330.6666
$ws.Range("L82").Value = 1355
$ws.Range("M82").Value = -969.6666
$ws.Range("N82").Value = -2077

$ws.Range("H85").Value = 1347.7
$ws.Range("I85").Value = 1330.6666
$ws.Range("J85").Value = 1355
$ws.Range("K85").Value = 1330.6666
$ws.Range("L85").Value = 1355
$ws.Range("M85").Value = -82.66660000000002
$ws.Range("N85").Value = -3851

$ws.Range("H136").Value = 3106.1738
$ws.Range("I136").Value = 727.2820400000001
$ws.Range("J136").Value = 16360
$ws.Range("K136").Value = 2181.84612
$ws.Range("L136").Value = 49080
$ws.Range("M136").Value = 368.1538799999998
$ws.Range("N136").Value = -54180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 6333.3335
$ws.Range("J40").Value = 6333.3335
$ws.Range("L40").Value = 6333.3335
$ws.Range("N40").Value = -6631.3335
